# Update the workbook for the "Add data for 2022-08-05" commit.
# This adds one more day (July 28) of carjacking data to the "July" columns
# for each year (2016-2022), and renames the sheet/header label accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet tab from "Through 2022-07-27" to "Through 2022-07-28"
$ws.Name = "Through 2022-07-28"

# 2. Update the header label for column B (shared string used in B1)
$ws.Range("B1").Value = "July 2022 (through July 28)"

# 3. Apply the updated / new cell values.
#    Mapping: row -> @{ column letter = new value }
$updates = @{
    "B2"  = 14
    "I2"  = 11
    "AD3" = 2
    "AK4" = 5
    "P6"  = 5
    "W8"  = 1
    "B11" = 3
    "I19" = 4
    "P19" = 2
    "B20" = 6
    "W21" = 1
    "AK23" = 3
    "AD26" = 1
    "B30" = 1
    "P33" = 1
    "AR38" = 3
    "AK39" = 4
    "W41" = 2
    "B45" = 3
    "AK48" = 1
    "B52" = 6
    "W62" = 2
    "I78" = 4
    "I82" = 1
    "P94" = 5
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
